$wb = $excel.ActiveWorkbook

# --- Rename Sheet2 -> InvalidLogin ---
$wsValid = $wb.Worksheets.Item("ValidLogin")
$wsInvalid = $wb.Worksheets.Item("Sheet2")
$wsInvalid.Name = "InvalidLogin"

# --- Populate InvalidLogin with header + invalid-credentials data row ---
$wsInvalid.Range("A1").Value = "UserName"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "abcd"
$wsInvalid.Range("B2").Value = "xyz"

# Reuse the same header/data cell formatting (bold+border header, bordered data)
# that is already used on the ValidLogin sheet, by copying formats across.
$wsValid.Range("A1:B1").Copy() | Out-Null
$wsInvalid.Range("A1:B1").PasteSpecial(-4122) | Out-Null

$wsValid.Range("A2:B2").Copy() | Out-Null
$wsInvalid.Range("A2:B2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Size the new columns to fit their content (UserName/abcd, Password/xyz).
$wsInvalid.Columns.Item(1).AutoFit() | Out-Null
$wsInvalid.Columns.Item(2).AutoFit() | Out-Null

# --- Selections: ValidLogin keeps a plain range selection, no longer the active tab ---
$wsValid.Range("A1:B2").Select() | Out-Null

# --- InvalidLogin becomes the active sheet/tab, selection resting on B2 ---
$wsInvalid.Activate() | Out-Null
$wsInvalid.Range("B2").Select() | Out-Null

# Match the zoom level captured for the InvalidLogin sheet view.
$excel.ActiveWindow.Zoom = 130
